$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 4036.5386
$ws_ALC.Range("I40").Value = 4884.1
$ws_ALC.Range("J40").Value = 1211.3334
$ws_ALC.Range("K40").Value = 4884.1
$ws_ALC.Range("L40").Value = 1211.3334
$ws_ALC.Range("M40").Value = -4709.1
$ws_ALC.Range("N40").Value = -1561.3334
$ws_ALC.Range("H64").Value = 3284.25
$ws_ALC.Range("I64").Value = 3597.3333
$ws_ALC.Range("J64").Value = 2345
$ws_ALC.Range("K64").Value = 3597.3333
$ws_ALC.Range("L64").Value = 2345
$ws_ALC.Range("M64").Value = -3349.3333
$ws_ALC.Range("N64").Value = -2841
$ws_ALC.Range("H67").Value = 3284.25
$ws_ALC.Range("I67").Value = 3597.3333
$ws_ALC.Range("J67").Value = 2345
$ws_ALC.Range("K67").Value = 3597.3333
$ws_ALC.Range("L67").Value = 2345
$ws_ALC.Range("M67").Value = -2739.3333
$ws_ALC.Range("N67").Value = -4061
$ws_ALC.Range("H74").Value = 6370.8
$ws_ALC.Range("I74").Value = 6370.8
$ws_ALC.Range("K74").Value = 6370.8
$ws_ALC.Range("M74").Value = -5434.8
$ws_ALC.Range("H77").Value = 6370.8
$ws_ALC.Range("I77").Value = 6370.8
$ws_ALC.Range("K77").Value = 31854
$ws_ALC.Range("M77").Value = -27174
$ws_ALC.Range("H80").Value = 997.5
$ws_ALC.Range("I80").Value = 548.625
$ws_ALC.Range("J80").Value = 1446.375
$ws_ALC.Range("K80").Value = 1645.875
$ws_ALC.Range("L80").Value = 4339.125
$ws_ALC.Range("M80").Value = -647.875
$ws_ALC.Range("N80").Value = -6335.125
$ws_ALC.Range("H83").Value = 997.5
$ws_ALC.Range("I83").Value = 548.625
$ws_ALC.Range("J83").Value = 1446.375
$ws_ALC.Range("K83").Value = 4937.625
$ws_ALC.Range("L83").Value = 13017.375
$ws_ALC.Range("M83").Value = 54.375
$ws_ALC.Range("N83").Value = -23001.375
$ws_ALC.Range("H98").Value = 314924.88
$ws_ALC.Range("I98").Value = 2269.8667
$ws_ALC.Range("J98").Value = 5004750
$ws_ALC.Range("K98").Value = 2269.8667
$ws_ALC.Range("L98").Value = 5004750
$ws_ALC.Range("M98").Value = -771.8667
$ws_ALC.Range("N98").Value = -5007746
$ws_ALC.Range("H122").Value = 314924.88
$ws_ALC.Range("I122").Value = 2269.8667
$ws_ALC.Range("J122").Value = 5004750
$ws_ALC.Range("K122").Value = 6809.6001
$ws_ALC.Range("L122").Value = 15014250
$ws_ALC.Range("M122").Value = -4359.6001
$ws_ALC.Range("N122").Value = -15019150
$ws_ALC.Range("H137").Value = 33976.39
$ws_ALC.Range("I137").Value = 81511.57000000001
$ws_ALC.Range("J137").Value = 3726.7273
$ws_ALC.Range("K137").Value = 244534.71
$ws_ALC.Range("L137").Value = 11180.1819
$ws_ALC.Range("M137").Value = -241984.71
$ws_ALC.Range("N137").Value = -16280.1819
$ws_ALC.Range("H138").Value = 2152.7778
$ws_ALC.Range("I138").Value = 1796.875
$ws_ALC.Range("J138").Value = 5000
$ws_ALC.Range("K138").Value = 5390.625
$ws_ALC.Range("L138").Value = 15000
$ws_ALC.Range("M138").Value = -250.625
$ws_ALC.Range("N138").Value = -25280

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 5018.5
$ws_ARM.Range("I32").Value = 2015.4048
$ws_ARM.Range("J32").Value = 30244.5
$ws_ARM.Range("K32").Value = 2015.4048
$ws_ARM.Range("L32").Value = 30244.5
$ws_ARM.Range("M32").Value = -1728.4048
$ws_ARM.Range("N32").Value = -30818.5
$ws_ARM.Range("H45").Value = 4870.25
$ws_ARM.Range("I45").Value = 4870.25
$ws_ARM.Range("K45").Value = 4870.25
$ws_ARM.Range("M45").Value = -4493.25
$ws_ARM.Range("H61").Value = 2401.1702
$ws_ARM.Range("I61").Value = 2225.513
$ws_ARM.Range("K61").Value = 2225.513
$ws_ARM.Range("M61").Value = -2013.513
$ws_ARM.Range("H63").Value = 3599
$ws_ARM.Range("I63").Value = 3749.25
$ws_ARM.Range("J63").Value = 2998
$ws_ARM.Range("K63").Value = 3749.25
$ws_ARM.Range("L63").Value = 2998
$ws_ARM.Range("M63").Value = -3063.25
$ws_ARM.Range("N63").Value = -4370
$ws_ARM.Range("H66").Value = 3599
$ws_ARM.Range("I66").Value = 3749.25
$ws_ARM.Range("J66").Value = 2998
$ws_ARM.Range("K66").Value = 18746.25
$ws_ARM.Range("L66").Value = 14990
$ws_ARM.Range("M66").Value = -15314.25
$ws_ARM.Range("N66").Value = -21854
$ws_ARM.Range("H132").Value = 2264.7144
$ws_ARM.Range("I132").Value = 2292.9092
$ws_ARM.Range("K132").Value = 6878.7276
$ws_ARM.Range("M132").Value = -4348.7276
$ws_ARM.Range("H136").Value = 2401.1702
$ws_ARM.Range("I136").Value = 2225.513
$ws_ARM.Range("K136").Value = 6676.539
$ws_ARM.Range("M136").Value = -4126.539

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 4208.1577
$ws_BSM.Range("I86").Value = 4089.5454
$ws_BSM.Range("K86").Value = 4089.5454
$ws_BSM.Range("M86").Value = -2966.5454
$ws_BSM.Range("H89").Value = 4208.1577
$ws_BSM.Range("I89").Value = 4089.5454
$ws_BSM.Range("K89").Value = 20447.727
$ws_BSM.Range("M89").Value = -14831.727
$ws_BSM.Range("H107").Value = 2150.182
$ws_BSM.Range("I107").Value = 1996.1111
$ws_BSM.Range("J107").Value = 2843.5
$ws_BSM.Range("K107").Value = 1996.1111
$ws_BSM.Range("L107").Value = 2843.5
$ws_BSM.Range("M107").Value = -76.11110000000008
$ws_BSM.Range("N107").Value = -6683.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H105").Value = 5497.512
$ws_CRP.Range("I105").Value = 1446.1
$ws_CRP.Range("J105").Value = 6804.4194
$ws_CRP.Range("K105").Value = 1446.1
$ws_CRP.Range("L105").Value = 6804.4194
$ws_CRP.Range("M105").Value = 300.9000000000001
$ws_CRP.Range("N105").Value = -10298.4194
$ws_CRP.Range("H132").Value = 5046.132
$ws_CRP.Range("I132").Value = 4672.927
$ws_CRP.Range("K132").Value = 14018.781
$ws_CRP.Range("M132").Value = -11488.781

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H68").Value = 2781556.2
$ws_CUL.Range("I68").Value = 8335551.5
$ws_CUL.Range("J68").Value = 4558.6
$ws_CUL.Range("K68").Value = 25006654.5
$ws_CUL.Range("L68").Value = 13675.8
$ws_CUL.Range("M68").Value = -25005843.5
$ws_CUL.Range("N68").Value = -15297.8
$ws_CUL.Range("H71").Value = 2781556.2
$ws_CUL.Range("I71").Value = 8335551.5
$ws_CUL.Range("J71").Value = 4558.6
$ws_CUL.Range("K71").Value = 75019963.5
$ws_CUL.Range("L71").Value = 41027.4
$ws_CUL.Range("M71").Value = -75015907.5
$ws_CUL.Range("N71").Value = -49139.4
$ws_CUL.Range("H107").Value = 509.36365
$ws_CUL.Range("I107").Value = 509.36365
$ws_CUL.Range("J107").Value = 0
$ws_CUL.Range("K107").Value = 1528.09095
$ws_CUL.Range("L107").Value = 0
$ws_CUL.Range("M107").Value = 391.90905
$ws_CUL.Range("N107").ClearContents()

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 3470.9375
$ws_GSM.Range("I70").Value = 2967.0908
$ws_GSM.Range("K70").Value = 2967.0908
$ws_GSM.Range("M70").Value = -2697.0908
$ws_GSM.Range("H73").Value = 3470.9375
$ws_GSM.Range("I73").Value = 2967.0908
$ws_GSM.Range("K73").Value = 2967.0908
$ws_GSM.Range("M73").Value = -2031.0908
$ws_GSM.Range("H113").Value = 2662.3333
$ws_GSM.Range("I113").Value = 2997
$ws_GSM.Range("J113").Value = 1993
$ws_GSM.Range("K113").Value = 2997
$ws_GSM.Range("L113").Value = 1993
$ws_GSM.Range("M113").Value = -827
$ws_GSM.Range("N113").Value = -6333
$ws_GSM.Range("H122").Value = 3363.5642
$ws_GSM.Range("I122").Value = 3217.5417
$ws_GSM.Range("J122").Value = 3597.2
$ws_GSM.Range("K122").Value = 9652.625100000001
$ws_GSM.Range("L122").Value = 10791.6
$ws_GSM.Range("M122").Value = -7202.625100000001
$ws_GSM.Range("N122").Value = -15691.6

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H69").Value = 12350
$ws_LTW.Range("I69").Value = 12350
$ws_LTW.Range("K69").Value = 12350
$ws_LTW.Range("M69").Value = -11539
$ws_LTW.Range("H72").Value = 12350
$ws_LTW.Range("I72").Value = 12350
$ws_LTW.Range("K72").Value = 37050
$ws_LTW.Range("M72").Value = -32994
$ws_LTW.Range("H122").Value = 4164.6787
$ws_LTW.Range("I122").Value = 4104.56
$ws_LTW.Range("K122").Value = 12313.68
$ws_LTW.Range("M122").Value = -9863.68

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 1029.9
$ws_WVR.Range("I107").Value = 1032.8889
$ws_WVR.Range("J107").Value = 1003
$ws_WVR.Range("K107").Value = 3098.6667
$ws_WVR.Range("L107").Value = 3009
$ws_WVR.Range("M107").Value = -1178.6667
$ws_WVR.Range("N107").Value = -6849
$ws_WVR.Range("H126").Value = 58075.652
$ws_WVR.Range("I126").Value = 2214.9524
$ws_WVR.Range("K126").Value = 6644.8572
$ws_WVR.Range("M126").Value = -4174.8572
$ws_WVR.Range("H136").Value = 358327.47
$ws_WVR.Range("I136").Value = 400966.97
$ws_WVR.Range("J136").Value = 2998.3333
$ws_WVR.Range("K136").Value = 1202900.91
$ws_WVR.Range("L136").Value = 8994.999899999999
$ws_WVR.Range("M136").Value = -1200350.91
$ws_WVR.Range("N136").Value = -14094.9999
